# Apply corrected quantity (F) and amount (G) values for affected stock rows,
# plus the resulting company Sub Total (B) and Grand Total (B1079/B1080) rollups.
# Amount (G) = Cost (D) * Quantity (F) for each line item; Sub/Grand totals are the
# sum of the Amount column for their respective block.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F137").Value = 7
$ws.Range("G137").Value = 792.8200000000001
$ws.Range("B139").Value = 92264.2
$ws.Range("F173").Value = 12
$ws.Range("G173").Value = 1721.28
$ws.Range("F180").Value = 58
$ws.Range("G180").Value = 5387.62
$ws.Range("F185").Value = 91
$ws.Range("G185").Value = 10481.38
$ws.Range("B186").Value = 57756
$ws.Range("F186").Value = 19
$ws.Range("G186").Value = 1262.36
$ws.Range("B187").Value = 53925
$ws.Range("F187").Value = 1
$ws.Range("G187").Value = 66.44
$ws.Range("F188").Value = 17
$ws.Range("G188").Value = 1645.77
$ws.Range("B191").Value = 79811.67999999999
$ws.Range("F256").Value = 7
$ws.Range("G256").Value = 463.61
$ws.Range("F259").Value = 26
$ws.Range("G259").Value = 2106
$ws.Range("B261").Value = 14338.05
$ws.Range("F272").Value = 139
$ws.Range("G272").Value = 4368.77
$ws.Range("F275").Value = 30
$ws.Range("G275").Value = 1209
$ws.Range("F282").Value = 30
$ws.Range("G282").Value = 2254.2
$ws.Range("F285").Value = 35
$ws.Range("G285").Value = 2468.55
$ws.Range("B293").Value = 61868.68
$ws.Range("F368").Value = 23
$ws.Range("G368").Value = 6802.25
$ws.Range("F372").Value = 58
$ws.Range("G372").Value = 4543.14
$ws.Range("F374").Value = 77
$ws.Range("G374").Value = 2469.39
$ws.Range("F395").Value = 23
$ws.Range("G395").Value = 1889.68
$ws.Range("B399").Value = 186262.06
$ws.Range("F409").Value = 167
$ws.Range("G409").Value = 28855.93
$ws.Range("F414").Value = 99
$ws.Range("G414").Value = 15929.1
$ws.Range("F422").Value = 41
$ws.Range("G422").Value = 5405.03
$ws.Range("F435").Value = 12
$ws.Range("G435").Value = 615.84
$ws.Range("F447").Value = 279
$ws.Range("G447").Value = 6450.48
$ws.Range("F449").Value = 140
$ws.Range("G449").Value = 15374.8
$ws.Range("F450").Value = 11
$ws.Range("G450").Value = 1172.05
$ws.Range("F458").Value = 49
$ws.Range("G458").Value = 4811.8
$ws.Range("F459").Value = 16
$ws.Range("G459").Value = 920.48
$ws.Range("F467").Value = 33
$ws.Range("G467").Value = 4105.86
$ws.Range("F469").Value = 557
$ws.Range("G469").Value = 32695.9
$ws.Range("F484").Value = 24
$ws.Range("G484").Value = 12622.8
$ws.Range("F485").Value = 480
$ws.Range("G485").Value = 82238.39999999999
$ws.Range("F497").Value = 288
$ws.Range("G497").Value = 17127.36
$ws.Range("F499").Value = 154
$ws.Range("G499").Value = 14265.02
$ws.Range("B503").Value = 446820.35
$ws.Range("F510").Value = 26
$ws.Range("G510").Value = 4204.98
$ws.Range("F513").Value = 11
$ws.Range("G513").Value = 2122.56
$ws.Range("B521").Value = 74657.03999999999
$ws.Range("F566").Value = 10
$ws.Range("G566").Value = 662.9
$ws.Range("F573").Value = 49
$ws.Range("G573").Value = 7087.36
$ws.Range("B575").Value = 89712.42
$ws.Range("F626").Value = 0
$ws.Range("G626").Value = 0
$ws.Range("F636").Value = 30
$ws.Range("G636").Value = 2159.1
$ws.Range("F644").Value = 4
$ws.Range("G644").Value = 143.56
$ws.Range("B645").Value = 25218.85
$ws.Range("F704").Value = 30
$ws.Range("G704").Value = 1388.4
$ws.Range("F705").Value = 164
$ws.Range("G705").Value = 5292.28
$ws.Range("F707").Value = 62
$ws.Range("G707").Value = 3001.42
$ws.Range("B715").Value = 127524.38
$ws.Range("F767").Value = 6
$ws.Range("G767").Value = 4507.8
$ws.Range("B772").Value = 7559.49
$ws.Range("F784").Value = 59
$ws.Range("G784").Value = 2483.31
$ws.Range("F786").Value = 114
$ws.Range("G786").Value = 30334.26
$ws.Range("F802").Value = 44
$ws.Range("G802").Value = 4242.48
$ws.Range("B803").Value = 142415.26
$ws.Range("F891").Value = 0
$ws.Range("G891").Value = 0
$ws.Range("B896").Value = 344710.76
$ws.Range("F918").Value = 21
$ws.Range("G918").Value = 3017.07
$ws.Range("F923").Value = 23
$ws.Range("G923").Value = 3017.14
$ws.Range("F927").Value = 23
$ws.Range("G927").Value = 3879.18
$ws.Range("F930").Value = 39
$ws.Range("G930").Value = 974.61
$ws.Range("F931").Value = 35
$ws.Range("G931").Value = 2175.95
$ws.Range("F934").Value = 27
$ws.Range("G934").Value = 4553.82
$ws.Range("B937").Value = 41570.27
$ws.Range("F954").Value = 155
$ws.Range("G954").Value = 12641.8
$ws.Range("F959").Value = 273
$ws.Range("G959").Value = 36336.3
$ws.Range("F964").Value = 127
$ws.Range("G964").Value = 4734.56
$ws.Range("F968").Value = 299
$ws.Range("G968").Value = 12393.55
$ws.Range("F969").Value = 173
$ws.Range("G969").Value = 9295.290000000001
$ws.Range("F970").Value = 439
$ws.Range("G970").Value = 66346.07000000001
$ws.Range("F974").Value = 144
$ws.Range("G974").Value = 20736
$ws.Range("B978").Value = 328826.95
$ws.Range("F989").Value = 197
$ws.Range("G989").Value = 29630.77
$ws.Range("F1002").Value = 383
$ws.Range("G1002").Value = 14105.89
$ws.Range("F1007").Value = 195
$ws.Range("G1007").Value = 9728.549999999999
$ws.Range("B1009").Value = 354302.75
$ws.Range("F1061").Value = 946
$ws.Range("G1061").Value = 154302.06
$ws.Range("F1063").Value = 72
$ws.Range("G1063").Value = 20366.64
$ws.Range("B1067").Value = 175166.14
$ws.Range("B1079").Value = 5109237.46
$ws.Range("B1080").Value = 5109237.46
